# Energy_share_requirement.xlsx cleanup
#
# The authors reorganised the example folder/module (see commit message)
# and, while doing so, re-saved the workbook from a fresh session:
#  - the stale "best fit" column widths on C:D were cleared back to the
#    sheet's normal/default column width, and
#  - the lingering multi-column selection (left over from the formatting
#    pass that produced the bold/centered header + wrapped description
#    cells) was replaced with a plain single-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C and D had a manual "best fit" width applied (22.83 chars) while
# the sheet was being laid out. Restore them to the sheet's normal column
# width instead of leaving that stale, oversized override in place.
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 8

# Reset the lingering selection (previously an entire C:D column range) to
# a simple single-cell selection.
$ws.Range("G2").Select()
